$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update extraction timestamps (kept as plain text, matching original formatting)
$ws.Range("C2").Value = "2025-01-25 19:07:18"
$ws.Range("C3").Value = "2025-01-25 17:10:47"
$ws.Range("C4").Value = "2025-01-25 16:31:01"
$ws.Range("C9").Value = "2025-01-25 16:19:09"

$ws.Range("C25").Value = "2025-01-25 17:25:41"
$ws.Range("D25").Value = "Macaé"
$ws.Range("E25").Value = "Rio de Janeiro"

$ws.Range("C32").Value = "2025-01-25 17:25:51"
$ws.Range("D32").Value = "Nova Iguaçu"
$ws.Range("E32").Value = "Rio de Janeiro"

$ws.Range("C53").Value = "2025-01-25 19:07:59"
$ws.Range("D53").Value = "Teresópolis"
$ws.Range("E53").Value = "Rio de Janeiro"

$ws.Range("C93").Value = "2025-01-25 17:10:54"
